# envio preliminar IRA 42 2025
# Adds week 42 data (column AS) to the weekly IRA extended report,
# updates a couple of institution names, and fills in a late AR28 value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header for week 42 (column AS), matching the style of the other
#     week-number header cells (bold + centered), but forced to text so it
#     stores like the existing "1".."41" inline-string headers. ---
$ws.Range("AS1").NumberFormat = "@"
$ws.Range("AS1").Value = "42"

# --- Institution name corrections (column C) ---
$ws.Range("C5").Value = "CAJA DE COMPENSACION FAMILIAR DE RISARALDA COMFAMI"
$ws.Range("C6").Value = "CAJA DE COMPENSACION FAMILIAR DE RISARALDA COMFAMI"
$ws.Range("C44").Value = "EMPRESA DE MEDICINA INTEGRAL EMI SA - SERVICIO DE"

# --- Late-arriving AR28 value (previously missing) ---
$ws.Range("AR28").Value = 210

# --- Week 42 counts (column AS) per institution row ---
$ws.Range("AS2").Value = 43
$ws.Range("AS3").Value = 46
$ws.Range("AS5").Value = 3
$ws.Range("AS7").Value = 19
$ws.Range("AS8").Value = 33
$ws.Range("AS9").Value = 1
$ws.Range("AS10").Value = 2
$ws.Range("AS11").Value = 3
$ws.Range("AS14").Value = 2
$ws.Range("AS16").Value = 3
$ws.Range("AS17").Value = 1
$ws.Range("AS21").Value = 1
$ws.Range("AS22").Value = 5
$ws.Range("AS23").Value = 7
$ws.Range("AS24").Value = 1
$ws.Range("AS25").Value = 24
$ws.Range("AS28").Value = 163
$ws.Range("AS29").Value = 0
$ws.Range("AS30").Value = 55
$ws.Range("AS31").Value = 3
$ws.Range("AS36").Value = 2
$ws.Range("AS37").Value = 12
$ws.Range("AS38").Value = 88
$ws.Range("AS41").Value = 10
$ws.Range("AS42").Value = 28
$ws.Range("AS43").Value = 24
$ws.Range("AS45").Value = 59
$ws.Range("AS46").Value = 75
$ws.Range("AS47").Value = 3
$ws.Range("AS48").Value = 38
$ws.Range("AS49").Value = 3
$ws.Range("AS50").Value = 0
$ws.Range("AS51").Value = 4
$ws.Range("AS53").Value = 9
$ws.Range("AS54").Value = 0
$ws.Range("AS55").Value = 2
$ws.Range("AS56").Value = 3
$ws.Range("AS57").Value = 16
$ws.Range("AS58").Value = 9
